# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, reflecting refreshed
# market-board price data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H33").Value = 110.42308
$ws.Range("I33").Value = 82.04348
$ws.Range("J33").Value = 328
$ws.Range("K33").Value = 82.04348
$ws.Range("L33").Value = 328
$ws.Range("M33").Value = 146.95652
$ws.Range("N33").Value = -786

$ws.Range("H107").Value = 890204.75
$ws.Range("I107").Value = 1270957.1
$ws.Range("J107").Value = 1782.6666
$ws.Range("K107").Value = 1270957.1
$ws.Range("L107").Value = 1782.6666
$ws.Range("M107").Value = -1269037.1
$ws.Range("N107").Value = -5622.6666

$ws.Range("H112").Value = 13637358
$ws.Range("J112").Value = 14355066
$ws.Range("L112").Value = 43065198
$ws.Range("N112").Value = -43067414

$ws.Range("H129").Value = 1017.5
$ws.Range("J129").Value = 1072.5778
$ws.Range("L129").Value = 3217.7334
$ws.Range("N129").Value = -13217.7334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2573.742
$ws.Range("I61").Value = 1698.2
$ws.Range("J61").Value = 4165.636
$ws.Range("K61").Value = 1698.2
$ws.Range("L61").Value = 4165.636
$ws.Range("M61").Value = -1486.2
$ws.Range("N61").Value = -4589.636

$ws.Range("H74").Value = 1274.8572
$ws.Range("I74").Value = 1458.8
$ws.Range("J74").Value = 1029.6
$ws.Range("K74").Value = 1458.8
$ws.Range("L74").Value = 1029.6
$ws.Range("M74").Value = -584.8
$ws.Range("N74").Value = -2777.6

$ws.Range("H77").Value = 1274.8572
$ws.Range("I77").Value = 1458.8
$ws.Range("J77").Value = 1029.6
$ws.Range("K77").Value = 7294
$ws.Range("L77").Value = 5148
$ws.Range("M77").Value = -2926
$ws.Range("N77").Value = -13884

$ws.Range("H132").Value = 3457.182
$ws.Range("I132").Value = 2870.9333
$ws.Range("J132").Value = 4713.4287
$ws.Range("K132").Value = 8612.7999
$ws.Range("L132").Value = 14140.2861
$ws.Range("M132").Value = -6082.7999
$ws.Range("N132").Value = -19200.2861

$ws.Range("H136").Value = 2573.742
$ws.Range("I136").Value = 1698.2
$ws.Range("J136").Value = 4165.636
$ws.Range("K136").Value = 5094.6
$ws.Range("L136").Value = 12496.908
$ws.Range("M136").Value = -2544.6
$ws.Range("N136").Value = -17596.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1870.3334
$ws.Range("I107").Value = 1625
$ws.Range("J107").Value = 2361
$ws.Range("K107").Value = 1625
$ws.Range("L107").Value = 2361
$ws.Range("M107").Value = 295
$ws.Range("N107").Value = -6201

$ws.Range("H134").Value = 2026.4407
$ws.Range("I134").Value = 1220.9796
$ws.Range("K134").Value = 3662.9388
$ws.Range("M134").Value = -1127.9388

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1193.26
$ws.Range("I31").Value = 907.59576
$ws.Range("J31").Value = 1446.585
$ws.Range("K31").Value = 907.59576
$ws.Range("L31").Value = 1446.585
$ws.Range("M31").Value = -612.59576
$ws.Range("N31").Value = -2036.585

$ws.Range("H34").Value = 1193.26
$ws.Range("I34").Value = 907.59576
$ws.Range("J34").Value = 1446.585
$ws.Range("K34").Value = 907.59576
$ws.Range("L34").Value = 1446.585
$ws.Range("M34").Value = -705.59576
$ws.Range("N34").Value = -1850.585

$ws.Range("H58").Value = 1660.4688
$ws.Range("I58").Value = 640.82355
$ws.Range("J58").Value = 2816.0667
$ws.Range("K58").Value = 640.82355
$ws.Range("L58").Value = 2816.0667
$ws.Range("M58").Value = -437.82355
$ws.Range("N58").Value = -3222.0667

$ws.Range("H107").Value = 533.88
$ws.Range("I107").Value = 445.93332
$ws.Range("K107").Value = 445.93332
$ws.Range("M107").Value = 1474.06668

$ws.Range("H122").Value = 906.8
$ws.Range("I122").Value = 783.5
$ws.Range("K122").Value = 2350.5
$ws.Range("M122").Value = 99.5

$ws.Range("H136").Value = 1660.4688
$ws.Range("I136").Value = 640.82355
$ws.Range("J136").Value = 2816.0667
$ws.Range("K136").Value = 1922.47065
$ws.Range("L136").Value = 8448.2001
$ws.Range("M136").Value = 627.5293500000002
$ws.Range("N136").Value = -13548.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 940.75
$ws.Range("I68").Value = 760.1746000000001
$ws.Range("J68").Value = 1248.2162
$ws.Range("K68").Value = 2280.5238
$ws.Range("L68").Value = 3744.6486
$ws.Range("M68").Value = -1469.5238
$ws.Range("N68").Value = -5366.6486

$ws.Range("H71").Value = 940.75
$ws.Range("I71").Value = 760.1746000000001
$ws.Range("J71").Value = 1248.2162
$ws.Range("K71").Value = 6841.571400000001
$ws.Range("L71").Value = 11233.9458
$ws.Range("M71").Value = -2785.571400000001
$ws.Range("N71").Value = -19345.9458

$ws.Range("H107").Value = 1062.1887
$ws.Range("I107").Value = 334.5
$ws.Range("K107").Value = 1003.5
$ws.Range("M107").Value = 916.5

$ws.Range("H122").Value = 744.9167
$ws.Range("I122").Value = 621.9375
$ws.Range("J122").Value = 843.3
$ws.Range("K122").Value = 5597.4375
$ws.Range("L122").Value = 7589.7
$ws.Range("M122").Value = -3147.4375
$ws.Range("N122").Value = -12489.7

$ws.Range("H131").Value = 2345.7979
$ws.Range("I131").Value = 499.0909
$ws.Range("J131").Value = 2606.2307
$ws.Range("K131").Value = 1497.2727
$ws.Range("L131").Value = 7818.6921
$ws.Range("M131").Value = 3542.7273
$ws.Range("N131").Value = -17898.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2146.5862
$ws.Range("I126").Value = 1912.4445
$ws.Range("J126").Value = 2251.95
$ws.Range("K126").Value = 5737.333500000001
$ws.Range("L126").Value = 6755.849999999999
$ws.Range("M126").Value = -3267.333500000001
$ws.Range("N126").Value = -11695.85

$ws.Range("H132").Value = 3835.2334
$ws.Range("I132").Value = 3485.238
$ws.Range("J132").Value = 4651.8887
$ws.Range("K132").Value = 10455.714
$ws.Range("L132").Value = 13955.6661
$ws.Range("M132").Value = -7925.714
$ws.Range("N132").Value = -19015.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 995
$ws.Range("I93").Value = 995
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 995
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 253

$ws.Range("H100").Value = 2235.879
$ws.Range("I100").Value = 1622.6666
$ws.Range("J100").Value = 2746.889
$ws.Range("K100").Value = 1622.6666
$ws.Range("L100").Value = 2746.889
$ws.Range("M100").Value = -1081.6666
$ws.Range("N100").Value = -3828.889

$ws.Range("H122").Value = 3746.8572
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 3983.3684
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 11950.1052
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -16850.1052

$ws.Range("H132").Value = 2919.2559
$ws.Range("I132").Value = 2001.3334
$ws.Range("J132").Value = 5037.5386
$ws.Range("K132").Value = 6004.0002
$ws.Range("L132").Value = 15112.6158
$ws.Range("M132").Value = -3474.0002
$ws.Range("N132").Value = -20172.6158

$ws.Range("H136").Value = 6068.5835
$ws.Range("I136").Value = 2275.3333
$ws.Range("J136").Value = 9861.833000000001
$ws.Range("K136").Value = 6825.999899999999
$ws.Range("L136").Value = 29585.499
$ws.Range("M136").Value = -4275.999899999999
$ws.Range("N136").Value = -34685.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 84270.914
$ws.Range("I126").Value = 100875.1
$ws.Range("J126").Value = 1250
$ws.Range("K126").Value = 302625.3
$ws.Range("L126").Value = 3750
$ws.Range("M126").Value = -300155.3
$ws.Range("N126").Value = -8690

$ws.Range("H136").Value = 22290780
$ws.Range("I136").Value = 37148860
$ws.Range("J136").Value = 3658.8333
$ws.Range("K136").Value = 111446580
$ws.Range("L136").Value = 10976.4999
$ws.Range("M136").Value = -111444030
$ws.Range("N136").Value = -16076.4999

